$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (90) of quotes, matching the style of the existing
# date column (style/format of A89) and the plain inline-string text
# cells used for the quote columns (B-E).
$row = 90

$ws.Cells.Item($row, 1).Value = 45995
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = "22,0574"
$ws.Cells.Item($row, 3).Value = "15,9382"
$ws.Cells.Item($row, 4).Value = "15,9382"
$ws.Cells.Item($row, 5).Value = "15,9382"
